$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) updates for rows 2-51.
# D-column values are written via NumberFormat "@" + Style reset so that
# numeric-looking price strings remain Text (matching the source data,
# and preserving trailing zeros / exact decimal text) instead of being
# auto-coerced into Excel Number values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.178.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.668.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5125"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.58%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2642"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06403"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07417"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.674.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.517"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5821"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008587"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.227.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.945"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.212"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.633"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1193"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.57%  "

$ws.Range("E27").Value = "  +2.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06416"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.298"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.320"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.67%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.532"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.527"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.640"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.015"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6097"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.366"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.655"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.162"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01606"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.083.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8659"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.70%  "

$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.816.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000112"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.100"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05206"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4293"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.905"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.21%  "
